# Update gh-pages output data for sheets "展览" (sheet 1) and "全部类型" (sheet 4).
# Both sheets contain the same table of events; column F ("想去人数") and
# column G ("最低票价") values are refreshed with newly scraped numbers.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 2;  F = 22 },
    @{ Row = 3;  F = 1814 },
    @{ Row = 4;  F = 35 },
    @{ Row = 5;  F = 782 },
    @{ Row = 11; F = 12 },
    @{ Row = 12; F = 149 },
    @{ Row = 15; F = 4234 },
    @{ Row = 18; F = 461 },
    @{ Row = 19; F = 400 },
    @{ Row = 20; F = 973 },
    @{ Row = 21; F = 1448 },
    @{ Row = 23; F = 40 },
    @{ Row = 25; F = 1960 },
    @{ Row = 27; F = 59; G = 50 },
    @{ Row = 28; F = 102 },
    @{ Row = 29; F = 196 }
)

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($update in $updates) {
        $row = $update.Row
        $ws.Cells.Item($row, 6).Value = $update.F
        if ($update.ContainsKey('G')) {
            $ws.Cells.Item($row, 7).Value = $update.G
        }
    }
}

$wb.Save()
